$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (2) through AC (29) hold the match data that gets
# swapped/rotated between rows; column A (the running index) stays put.
$firstCol = 2
$lastCol = 29

function Get-RowValues($row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += ,$ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$i]
        $i++
    }
}

# --- Swap rows 16 and 17 (all columns B:AC) ---
$row16 = Get-RowValues 16
$row17 = Get-RowValues 17
Set-RowValues 16 $row17
Set-RowValues 17 $row16

# --- Rotate rows 98, 99, 100 (all columns B:AC) ---
# new98 = old99 ; new99 = old100 ; new100 = old98
$row98 = Get-RowValues 98
$row99 = Get-RowValues 99
$row100 = Get-RowValues 100
Set-RowValues 98 $row99
Set-RowValues 99 $row100
Set-RowValues 100 $row98

# --- Delete row 123 (last data row) entirely ---
$ws.Rows.Item(123).Delete()
